# Updates cryptos list values (Price / Volume(1h) columns) per source refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.801.20'
$ws.Range('D3').Value = '1.632.16'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5014'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('E7').Value = '  -0.48%  '
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.64'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07700'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').Value = '1.664.80'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.249'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').Value = '1.857.20'
$ws.Range('E14').Value = '  -1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5413'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.81%  '
$ws.Range('D16').Value = '0.0₅7904'
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.51'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '25.810.97'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  -0.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '200.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.322'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.876'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.928'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.932'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +11.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.73%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1134'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.63'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.691'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.06%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04982'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.264'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.184'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.05%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '1.166.24'
$ws.Range('E36').Value = '  +0.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.623'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8892'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5560'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('E40').Value = '  -2.06%  '
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.667'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8037'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').Value = '1.769.80'
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4511'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.000'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '54.56'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05070'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.45%  '

Write-Host "Applied cryptos update"
